# The deck's Design theme ("Integral") is replaced with the stock default
# "Office Theme" colour palette (the same swap PowerPoint performs when a
# new theme is applied from the Design gallery).
#
# A PowerPoint theme's 12 colour scheme slots are exposed on the object
# model as ThemeColorScheme.Colors(1..12), in the fixed order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2,
#   5 accent1, 6 accent2, 7 accent3, 8 accent4, 9 accent5, 10 accent6,
#   11 hyperlink, 12 followed hyperlink
# Each ThemeColor's .RGB is the usual Win32 RGB() packing (R + G*256 + B*65536).

function Hex2Rgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# The presentation's overall design is carried by the slide master's theme.
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

$officeColors = @{
    1  = "000000" # dk1
    2  = "FFFFFF" # lt1
    3  = "44546A" # dk2
    4  = "E7E6E6" # lt2
    5  = "5B9BD5" # accent1
    6  = "ED7D31" # accent2
    7  = "A5A5A5" # accent3
    8  = "FFC000" # accent4
    9  = "4472C4" # accent5
    10 = "70AD47" # accent6
    11 = "0563C1" # hyperlink
    12 = "954F72" # followed hyperlink
}

foreach ($slot in 1..12) {
    $colorScheme.Colors($slot).RGB = Hex2Rgb $officeColors[$slot]
}
